$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 112 (pushes existing rows 112..190 down to 113..191),
# matching the row format of the row above it (date style carries over).
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the new daily price record.
$ws.Range("A112").Value = 5
$ws.Range("B112").Value = "Macroferia Regional de Talca"
$ws.Range("C112").Value = "Maule"
$ws.Range("D112").Value = 45236
$ws.Range("E112").Value = 7
$ws.Range("F112").Value = "Fruta"
$ws.Range("G112").Value = 100107
$ws.Range("H112").Value = "Otros"
$ws.Range("I112").Value = 100107002
$ws.Range("J112").Value = "Chirimoya"
$ws.Range("K112").Value = "Cultivar IV Región"
$ws.Range("L112").Value = "Segunda"
$ws.Range("M112").Value = 250
$ws.Range("N112").Value = 18000
$ws.Range("O112").Value = 18000
$ws.Range("P112").Value = 18000
$ws.Range("Q112").Value = "$/bandeja 10 kilos"
$ws.Range("R112").Value = "Provincia de Limarí"
$ws.Range("S112").Value = 1800
$ws.Range("T112").Value = 10
